# Diagram "exercise_tools" slide cleanup:
#  - Remove the standalone "Diagram exercise_tools" title textbox.
#  - Group the remaining picture, the two callout textboxes, and the two
#    connector arrows into a single "Group 1" group shape (matches
#    selecting the shapes in the UI and pressing Ctrl+G).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate shapes by name so the script is robust to any pre-existing
# selection/order assumptions.
$titleBox = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 24") {
        $titleBox = $s.Shapes.Item($i)
    }
}
if ($titleBox -ne $null) {
    $titleBox.Delete()
}

# Collect the indices of the shapes that remain (picture, the two
# callout textboxes, and the two connector arrows) and group them.
$idxList = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $idxList += $i
}

$range = $s.Shapes.Range($idxList)
$grp = $range.Group()
$grp.Name = "Group 1"
